$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before the old "issueDatepublished" column (old F),
#     shifting every following column one to the right. Excel's native
#     "insert column" naturally carries the left-neighbour's column width
#     onto the new column and keeps the bold header style on the new header
#     cell. ---
$ws.Columns("F:F").Insert()

# The new column inherits the same width as its left neighbours (D:E).
$ws.Columns("F:F").ColumnWidth = $ws.Columns("D:D").ColumnWidth

# New header cell for the inserted column.
$ws.Range("F1").Value = "language"

# --- New row 5: Finnish-language variant of the "New Issue, new article"
#     entry in row 4 (alternative-locale support). ---
$ws.Range("F5").Value = "fi"
$ws.Range("B5").Value = "Artikkeli suomeksi"
$ws.Range("D5").Value = "Tämä artikkeli on suomeksi"

$ws.Range("G5").Value = 42885
$ws.Range("G5").NumberFormat = "yyyy\-mm\-dd;@"

$ws.Range("H5").Value = 29
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 2017

$ws.Range("L5").Value = "Articles"
$ws.Range("M5").Value = "ART"

$ws.Range("N5").Value = "35-45"
$ws.Range("N5").NumberFormat = "@"
$ws.Range("O5").Value = "2"
$ws.Range("O5").NumberFormat = "@"

$ws.Range("P5").Value = "Antti-Jussi"
$ws.Range("Q5").Value = "Nygård"
$ws.Range("R5").Value = "Tieteellisten seurain valtuuskunta"

$ws.Range("V5").Value = "2017-2-2-1.pdf"
$ws.Range("W5").Value = "PDF"
$ws.Range("X5").Value = "Article Text"
$ws.Range("Y5").Value = "fi"

# --- View state: scroll/selection matching the author's saved session. ---
$ws.Range("Z5").Select()
